$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 500.97058
$ws.Range("I15").Value = 500.97058
$ws.Range("K15").Value = 1502.91174
$ws.Range("M15").Value = -1333.91174
$ws.Range("H28").Value = 2861.25
$ws.Range("I28").Value = 2315
$ws.Range("J28").Value = 4500
$ws.Range("K28").Value = 2315
$ws.Range("L28").Value = 4500
$ws.Range("M28").Value = -1830
$ws.Range("N28").Value = -5470
$ws.Range("H29").Value = 3130.8
$ws.Range("J29").Value = 6002
$ws.Range("L29").Value = 18006
$ws.Range("N29").Value = -18568
$ws.Range("H51").Value = 3344162.8
$ws.Range("J51").Value = 7187.25
$ws.Range("L51").Value = 7187.25
$ws.Range("N51").Value = -8155.25
$ws.Range("H76").Value = 4531
$ws.Range("I76").Value = 4236.6665
$ws.Range("J76").Value = 4972.5
$ws.Range("K76").Value = 4236.6665
$ws.Range("L76").Value = 4972.5
$ws.Range("M76").Value = -3921.6665
$ws.Range("N76").Value = -5602.5
$ws.Range("H79").Value = 4531
$ws.Range("I79").Value = 4236.6665
$ws.Range("J79").Value = 4972.5
$ws.Range("K79").Value = 4236.6665
$ws.Range("L79").Value = 4972.5
$ws.Range("M79").Value = -3144.6665
$ws.Range("N79").Value = -7156.5
$ws.Range("H107").Value = 1250.5625
$ws.Range("I107").Value = 1358.4615
$ws.Range("K107").Value = 1358.4615
$ws.Range("M107").Value = 561.5385000000001
$ws.Range("H111").Value = 17445.363
$ws.Range("I111").Value = 28666.666
$ws.Range("K111").Value = 85999.99800000001
$ws.Range("M111").Value = -82932.99800000001
$ws.Range("H127").Value = 1588.3334
$ws.Range("I127").Value = 879.4
$ws.Range("K127").Value = 2638.2
$ws.Range("M127").Value = 2321.8
$ws.Range("H138").Value = 3093.5
$ws.Range("J138").Value = 8255.286
$ws.Range("L138").Value = 24765.858
$ws.Range("N138").Value = -35045.858

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4179.2144
$ws.Range("I2").Value = 4278.778
$ws.Range("K2").Value = 4278.778
$ws.Range("M2").Value = -4165.778
$ws.Range("H32").Value = 1948.127
$ws.Range("I32").Value = 1737.2222
$ws.Range("J32").Value = 3213.5557
$ws.Range("K32").Value = 1737.2222
$ws.Range("L32").Value = 3213.5557
$ws.Range("M32").Value = -1450.2222
$ws.Range("N32").Value = -3787.5557
$ws.Range("H61").Value = 4767.8423
$ws.Range("I61").Value = 2839.7568
$ws.Range("K61").Value = 2839.7568
$ws.Range("M61").Value = -2627.7568
$ws.Range("H74").Value = 1360.0698
$ws.Range("I74").Value = 642.4286
$ws.Range("K74").Value = 642.4286
$ws.Range("M74").Value = 231.5714
$ws.Range("H77").Value = 1360.0698
$ws.Range("I77").Value = 642.4286
$ws.Range("K77").Value = 3212.143
$ws.Range("M77").Value = 1155.857
$ws.Range("H97").Value = 6968.3335
$ws.Range("I97").Value = 1772
$ws.Range("J97").Value = 17361
$ws.Range("K97").Value = 1772
$ws.Range("L97").Value = 17361
$ws.Range("M97").Value = -1276
$ws.Range("N97").Value = -18353
$ws.Range("H114").Value = 69999
$ws.Range("J114").Value = 69999
$ws.Range("L114").Value = 69999
$ws.Range("N114").Value = -78677
$ws.Range("H116").Value = 4179.2144
$ws.Range("I116").Value = 4278.778
$ws.Range("K116").Value = 4278.778
$ws.Range("M116").Value = -1984.778
$ws.Range("H122").Value = 2940.6572
$ws.Range("I122").Value = 2952.8215
$ws.Range("J122").Value = 2892
$ws.Range("K122").Value = 8858.4645
$ws.Range("L122").Value = 8676
$ws.Range("M122").Value = -6408.4645
$ws.Range("N122").Value = -13576
$ws.Range("H136").Value = 4767.8423
$ws.Range("I136").Value = 2839.7568
$ws.Range("K136").Value = 8519.270400000001
$ws.Range("M136").Value = -5969.270400000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4179.2144
$ws.Range("I3").Value = 4278.778
$ws.Range("K3").Value = 4278.778
$ws.Range("M3").Value = -4164.778

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2508.0715
$ws.Range("J16").Value = 3372.25
$ws.Range("L16").Value = 3372.25
$ws.Range("N16").Value = -3946.25
$ws.Range("H31").Value = 2032.3636
$ws.Range("I31").Value = 1144.7354
$ws.Range("J31").Value = 3469.476
$ws.Range("K31").Value = 1144.7354
$ws.Range("L31").Value = 3469.476
$ws.Range("M31").Value = -849.7354
$ws.Range("N31").Value = -4059.476
$ws.Range("H34").Value = 2032.3636
$ws.Range("I34").Value = 1144.7354
$ws.Range("J34").Value = 3469.476
$ws.Range("K34").Value = 1144.7354
$ws.Range("L34").Value = 3469.476
$ws.Range("M34").Value = -942.7354
$ws.Range("N34").Value = -3873.476
$ws.Range("H41").Value = 8530.5
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10856
$ws.Range("H62").Value = 3853.7144
$ws.Range("I62").Value = 3746
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 3746
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -3122
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 3853.7144
$ws.Range("I65").Value = 3746
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 18730
$ws.Range("L65").Value = 4500
$ws.Range("M65").Value = -15610
$ws.Range("N65").Value = -28740
$ws.Range("H99").Value = 6213.5654
$ws.Range("I99").Value = 6003.15
$ws.Range("K99").Value = 6003.15
$ws.Range("M99").Value = -4505.15
$ws.Range("H107").Value = 932.875
$ws.Range("I107").Value = 899.2
$ws.Range("K107").Value = 899.2
$ws.Range("M107").Value = 1020.8
$ws.Range("H113").Value = 2508.0715
$ws.Range("J113").Value = 3372.25
$ws.Range("L113").Value = 3372.25
$ws.Range("N113").Value = -7712.25
$ws.Range("H126").Value = 6213.5654
$ws.Range("I126").Value = 6003.15
$ws.Range("K126").Value = 18009.45
$ws.Range("M126").Value = -15539.45

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 8815.154
$ws.Range("J38").Value = 13987.375
$ws.Range("L38").Value = 41962.125
$ws.Range("N38").Value = -42656.125
$ws.Range("H55").Value = 8656.714
$ws.Range("J55").Value = 9916.5
$ws.Range("L55").Value = 29749.5
$ws.Range("N55").Value = -30103.5
$ws.Range("H60").Value = 955.8570999999999
$ws.Range("I60").Value = 322.75
$ws.Range("K60").Value = 968.25
$ws.Range("M60").Value = -717.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6094.5835
$ws.Range("I80").Value = 4017.6
$ws.Range("K80").Value = 4017.6
$ws.Range("M80").Value = -3019.6
$ws.Range("H82").Value = 80298
$ws.Range("I82").Value = 80298
$ws.Range("K82").Value = 80298
$ws.Range("M82").Value = -79915
$ws.Range("H83").Value = 6094.5835
$ws.Range("I83").Value = 4017.6
$ws.Range("K83").Value = 20088
$ws.Range("M83").Value = -15096
$ws.Range("H85").Value = 80298
$ws.Range("I85").Value = 80298
$ws.Range("K85").Value = 80298
$ws.Range("M85").Value = -78972
$ws.Range("H102").Value = 27835.45
$ws.Range("I102").Value = 2247.4666
$ws.Range("K102").Value = 2247.4666
$ws.Range("M102").Value = -625.4666000000002
$ws.Range("H113").Value = 3277
$ws.Range("I113").Value = 4124.25
$ws.Range("K113").Value = 4124.25
$ws.Range("M113").Value = -1954.25
$ws.Range("H126").Value = 23946.9
$ws.Range("I126").Value = 3308.625
$ws.Range("J126").Value = 106500
$ws.Range("K126").Value = 9925.875
$ws.Range("L126").Value = 319500
$ws.Range("M126").Value = -7455.875
$ws.Range("N126").Value = -324440

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -5272
$ws.Range("H46").Value = 6476.6855
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 6999.5
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 6999.5
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -7375.5
$ws.Range("H104").Value = 43318.89
$ws.Range("J104").Value = 43318.89
$ws.Range("L104").Value = 43318.89
$ws.Range("N104").Value = -50306.89
$ws.Range("H122").Value = 7009.657
$ws.Range("I122").Value = 6385.5938
$ws.Range("J122").Value = 13666.333
$ws.Range("K122").Value = 19156.7814
$ws.Range("L122").Value = 40998.999
$ws.Range("M122").Value = -16706.7814
$ws.Range("N122").Value = -45898.999
$ws.Range("H136").Value = 4935.9165
$ws.Range("I136").Value = 4930.1514
$ws.Range("J136").Value = 4999.3335
$ws.Range("K136").Value = 14790.4542
$ws.Range("L136").Value = 14998.0005
$ws.Range("M136").Value = -12240.4542
$ws.Range("N136").Value = -20098.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 880.5625
$ws.Range("I107").Value = 737.25
$ws.Range("K107").Value = 2211.75
$ws.Range("M107").Value = -291.75
$ws.Range("H122").Value = 2016.9318
$ws.Range("I122").Value = 2016.9318
$ws.Range("K122").Value = 6050.7954
$ws.Range("M122").Value = -3600.7954
$ws.Range("H136").Value = 4503.3726
$ws.Range("I136").Value = 4281.341
$ws.Range("K136").Value = 12844.023
$ws.Range("M136").Value = -10294.023
